$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data rows got shuffled: the Fecha/Volumen/Precio columns (D, M, N, O, P, S)
# for rows 2,3,6,7,8,9,10 now carry the values that used to belong to a different row.
# Row 4 and Row 5 are unchanged.

# Target values per row (Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio, Precio $/Kg)
$rowData = @{
    2  = @(44798, 80,  21000, 22000, 21500, 1075)
    3  = @(44357, 100, 14000, 15000, 14500, 725)
    6  = @(44761, 100, 20000, 21000, 20500, 1025)
    7  = @(44893, 80,  21000, 22000, 21625, 1081)
    8  = @(44792, 100, 21000, 22000, 21500, 1075)
    9  = @(44320, 80,  16000, 17000, 16500, 825)
    10 = @(44708, 80,  20000, 21000, 20500, 1025)
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    $ws.Range("D$row").Value = $vals[0]
    $ws.Range("M$row").Value = $vals[1]
    $ws.Range("N$row").Value = $vals[2]
    $ws.Range("O$row").Value = $vals[3]
    $ws.Range("P$row").Value = $vals[4]
    $ws.Range("S$row").Value = $vals[5]
}
